# Rename "form_id" setting to "table_id" in the settings sheet,
# and add a new "properties" sheet summarizing table/partition/aspect/key data,
# as described in the commit message:
#  "Update to some ODK Survey forms to rename them to the table_id so that we
#   generate definitions.csv and properties.csv ; update to process a
#   properties sheet into the properties.csv ; minimize the content of the
#   properties.csv"

$wb = $excel.ActiveWorkbook

# --- 1. settings sheet: rename "form_id" setting name to "table_id" ---
$settings = $wb.Worksheets.Item("settings")
$settings.Range("A2").Value = "table_id"

# --- 2. add a new "properties" sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$props = $wb.Worksheets.Add($null, $lastSheet)
$props.Name = "properties"

$props.Range("A1").Value = "partition"
$props.Range("B1").Value = "aspect"
$props.Range("C1").Value = "key"
$props.Range("D1").Value = "type"
$props.Range("E1").Value = "value"

$props.Range("A2").Value = "Table"
$props.Range("B2").Value = "default"
$props.Range("C2").Value = "colOrder"
$props.Range("D2").Value = "array"
$props.Range("E2").Value = '["FMT_FOL_date","FMT_FOL_B_focal_AnimID","FMT_time","FMT_FMP_seq_num","FMT_notes","FMT_fixed","FMT_FOL_foll_num"]'

# column widths to roughly match the authored sheet (~13.9 on A:E)
$props.Range("A1:E1").ColumnWidth = 13

# --- 3. update the "settings" sheet's selection (A3) without leaving it active ---
$settings.Range("A3").Select()

# --- 4. make "properties" the active/selected sheet, with E5 selected ---
$props.Activate()
$props.Range("E5").Select()
